$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix capitalization: "Paid Work" -> "Paid work" in the figure's category label (A8)
$ws.Range("A8").Value = "Paid work"

# Update the saved selection to A8 (matches author's cursor position when saving)
$ws.Range("A8").Select()
